$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values for rows 3-7 across the affected columns (A,B,E,F,G,H,Q,R)
$cols = @("A","B","E","F","G","H","Q","R")
$rows = @(3,4,5,6,7)

$data = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $data[$r] = $rowData
}

# New row order: target row N gets the data that was previously in row (mapping)
# target3 <- current6, target4 <- current7, target5 <- current3, target6 <- current4, target7 <- current5
$mapping = @{3=6; 4=7; 5=3; 6=4; 7=5}

foreach ($t in $rows) {
    $src = $mapping[$t]
    foreach ($c in $cols) {
        $ws.Range("$c$t").Value2 = $data[$src][$c]
    }
}
